$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.853.54'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '1.894.86'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7905'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.48'
$ws.Range("E6").Value = '  +0.55%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3148'
$ws.Range("E8").Value = '  -4.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.36'
$ws.Range("E9").Value = '  -4.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07260'
$ws.Range("E10").Value = '  +2.94%  '
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.575'
$ws.Range("E12").Value = '  +5.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7651'
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("D14").Value = '1.909.37'
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.47'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.174'
$ws.Range("E16").Value = '  +5.22%  '
$ws.Range("D17").Value = '29.859.51'
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.00'
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007790'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.175'
$ws.Range("E22").Value = '  +16.76%  '
$ws.Range("D23").Value = '2.152.77'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1655'
$ws.Range("E25").Value = '  -5.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.405'
$ws.Range("E26").Value = '  +1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.78'
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.061'
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.400'
$ws.Range("E30").Value = '  +2.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.545'
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.471'
$ws.Range("E32").Value = '  +4.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.092'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05540'
$ws.Range("E34").Value = '  -6.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.271'
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7395'
$ws.Range("E36").Value = '  +0.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9990'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.618'
$ws.Range("E38").Value = '  -3.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01924'
$ws.Range("E39").Value = '  +0.24%  '
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '1.145.15'
$ws.Range("E41").Value = '  +14.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.24'
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4421'
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.880'
$ws.Range("E44").Value = '  -1.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8504'
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.26'
$ws.Range("E46").Value = '  +2.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.878'
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.00'
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.449'
$ws.Range("E50").Value = '  -1.49%  '
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.009'
$ws.Range("E51").Value = '  +10.84%  '
